# feat: add robot to the project
# Update the "Webshop" sheet: fix Currys' URL, add two new webshop rows
# (hobbycraft, hm), and leave the Webshop tab as the active/selected sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Webshop")

# Fix the Currys URL (was a typo: "Curryes.co.uk")
$ws.Range("B4").Value = "https://www.currys.co.uk/"

# Add the new webshops
$ws.Range("A5").Value = "hobbycraft"
$ws.Range("B5").Value = "https://www.hobbycraft.co.uk/"

$ws.Range("A6").Value = "hm"
$ws.Range("B6").Value = "hm.com"

# Match formatting applied to the new URL cells
$ws.Range("B5:B6").Style = "Normal"

# Resize column B to fit the new (longer) URL text
$ws.Columns.Item(2).AutoFit()

# Make Webshop the active sheet / tab, with the next empty row selected
$ws.Activate()
$ws.Range("A7").Select()
